$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "_GoBack" bookmark that previously sat between the
#    "Armor 2/4" text and the trailing "/0" run.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Armor 2/4/0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="37F9E813" w14:textId="5E688D0F" w:rsidR="00256F46" w:rsidRDefault="00256F46" w:rsidP="00256F46"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:line="276" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t xml:space="preserve">Armor </w:t></w:r><w:r w:rsidR="00C13EE9"><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t>/</w:t></w:r><w:r w:rsidR="00C13EE9"><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t>4</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t>/0</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

# ------------------------------------------------------------------
# 2) Split "2d8 + Power energy damage" into four separate runs reading
#    "2d8 + Power " / "Energy" / " " / "Damage" / " (Reflex)", and move
#    the "_GoBack" bookmark to the end of this paragraph.
# ------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("2d8 + Power energy damage", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng2.Find.Found) {
    $xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="46901599" w14:textId="5685E340" w:rsidR="00256F46" w:rsidRDefault="00256F46" w:rsidP="00256F46"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:line="276" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t xml:space="preserve">2d8 + Power </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t>Energy</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t>Damage</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t xml:space="preserve"> (Reflex)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng2.InsertXML($xml2)
}

# ------------------------------------------------------------------
# 3) Delete the whole "Save for ½ damage (Reflex)" bullet paragraph.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Save for*damage (Reflex)*") {
        $p.Range.Delete()
        break
    }
}
